# Make MAIR plots (parity, histogram)
# The analysis notebook recomputed the "PerformerOverpassID" (column G)
# numbering after dropping two overpasses (rows 2 and 11) from the
# sequence, so every subsequent ID shifted down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows whose PerformerOverpassID (column G) cell is removed entirely.
$clearRows = @(2, 11)
foreach ($r in $clearRows) {
    $ws.Range("G$r").ClearContents() | Out-Null
}

# Remaining rows get renumbered PerformerOverpassID values.
$newValues = @{
    3  = 1
    4  = 2
    5  = 3
    6  = 4
    7  = 5
    8  = 6
    9  = 7
    10 = 8
    12 = 9
    13 = 10
    14 = 11
    15 = 12
    16 = 13
    17 = 14
    18 = 15
    19 = 16
    20 = 17
    21 = 18
    22 = 19
    23 = 20
    24 = 21
    25 = 22
    26 = 23
    27 = 24
}

foreach ($r in $newValues.Keys) {
    $ws.Range("G$r").Value = $newValues[$r]
}

# Column A (Date) was auto-fit to its content width while reviewing the
# new plots.
$ws.Columns.Item(1).AutoFit() | Out-Null
